{"js": "// Replace each two-digit-division answer cell's old text with the new text.\n// Each (old -> new) pair corresponds to one table cell, processed in document order.\n// Using Range.search + insertText(..., replace) preserves the run/paragraph formatting\n// (font, size, alignment) already present on the matched text.\nconst replacements = [[\"12\u00f75=2, 2\", \"78\u00f78=9, 6\"], [\"42\u00f76=7, 0\", \"52\u00f75=10, 2\"], [\"82\u00f74=20, 2\", \"68\u00f75=13, 3\"], [\"98\u00f78=12, 2\", \"85\u00f74=21, 1\"], [\"25\u00f77=3, 4\", \"31\u00f78=3, 7\"], [\"78\u00f76=13, 0\", \"40\u00f79=4, 4\"], [\"57\u00f76=9, 3\", \"21\u00f76=3, 3\"], [\"28\u00f79=3, 1\", \"44\u00f78=5, 4\"], [\"34\u00f77=4, 6\", \"41\u00f72=20, 1\"], [\"53\u00f77=7, 4\", \"54\u00f73=18, 0\"], [\"27\u00f76=4, 3\", \"82\u00f74=20, 2\"], [\"93\u00f77=13, 2\", \"82\u00f79=9, 1\"], [\"55\u00f78=6, 7\", \"82\u00f73=27, 1\"], [\"76\u00f73=25, 1\", \"93\u00f72=46, 1\"], [\"66\u00f76=11, 0\", \"78\u00f77=11, 1\"], [\"32\u00f73=10, 2\", \"19\u00f78=2, 3\"], [\"28\u00f73=9, 1\", \"38\u00f79=4, 2\"], [\"45\u00f76=7, 3\", \"70\u00f74=17, 2\"], [\"54\u00f74=13, 2\", \"35\u00f77=5, 0\"], [\"54\u00f77=7, 5\", \"17\u00f79=1, 8\"], [\"45\u00f79=5, 0\", \"53\u00f74=13, 1\"], [\"84\u00f72=42, 0\", \"54\u00f74=13, 2\"], [\"29\u00f74=7, 1\", \"76\u00f73=25, 1\"], [\"75\u00f75=15, 0\", \"96\u00f74=24, 0\"], [\"51\u00f72=25, 1\", \"99\u00f74=24, 3\"]];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  // Replace just the first (and expected only) match for this exact old value.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-division answer cell's old text with the new text.\n# Each (old -> new) pair corresponds to one table cell, processed in document order.\n# Using Find/Replace on the document Range preserves the run/paragraph formatting\n# (font, size, alignment) already present on the matched text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('12\u00f75=2, 2', '78\u00f78=9, 6'),\n    @('42\u00f76=7, 0', '52\u00f75=10, 2'),\n    @('82\u00f74=20, 2', '68\u00f75=13, 3'),\n    @('98\u00f78=12, 2', '85\u00f74=21, 1'),\n    @('25\u00f77=3, 4', '31\u00f78=3, 7'),\n    @('78\u00f76=13, 0', '40\u00f79=4, 4'),\n    @('57\u00f76=9, 3', '21\u00f76=3, 3'),\n    @('28\u00f79=3, 1', '44\u00f78=5, 4'),\n    @('34\u00f77=4, 6', '41\u00f72=20, 1'),\n    @('53\u00f77=7, 4', '54\u00f73=18, 0'),\n    @('27\u00f76=4, 3', '82\u00f74=20, 2'),\n    @('93\u00f77=13, 2', '82\u00f79=9, 1'),\n    @('55\u00f78=6, 7', '82\u00f73=27, 1'),\n    @('76\u00f73=25, 1', '93\u00f72=46, 1'),\n    @('66\u00f76=11, 0', '78\u00f77=11, 1'),\n    @('32\u00f73=10, 2', '19\u00f78=2, 3'),\n    @('28\u00f73=9, 1', '38\u00f79=4, 2'),\n    @('45\u00f76=7, 3', '70\u00f74=17, 2'),\n    @('54\u00f74=13, 2', '35\u00f77=5, 0'),\n    @('54\u00f77=7, 5', '17\u00f79=1, 8'),\n    @('45\u00f79=5, 0', '53\u00f74=13, 1'),\n    @('84\u00f72=42, 0', '54\u00f74=13, 2'),\n    @('29\u00f74=7, 1', '76\u00f73=25, 1'),\n    @('75\u00f75=15, 0', '96\u00f74=24, 0'),\n    @('51\u00f72=25, 1', '99\u00f74=24, 3'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdReplaceOne = 1 -> replace just the single (expected unique) match for this cell's text.\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
